$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Автобус 2 (06:00 - 14:00)"
$ws.Range("D2").Value = "Автобус 6 (07:00 - 15:00)"
$ws.Range("E2").Value = "Автобус 5 (09:00 - 17:00)"
$ws.Range("F2").Value = "Автобус 2 (07:00 - 15:00)"
$ws.Range("G2").Value = "Автобус 8 (08:00 - 16:00)"

$ws.Range("C3").Value = "Автобус 1 (09:00 - 17:00)"
$ws.Range("D3").Value = "Автобус 7 (08:00 - 16:00)"
$ws.Range("E3").Value = "Автобус 2 (10:00 - 18:00)"
$ws.Range("F3").Value = "Автобус 8 (08:00 - 16:00)"
$ws.Range("G3").Value = "Автобус 7 (08:00 - 16:00)"

$ws.Range("C4").Value = "Автобус 6 (08:00 - 16:00)"
$ws.Range("D4").Value = "Автобус 4 (09:00 - 17:00)"
$ws.Range("E4").Value = "Автобус 6 (07:00 - 15:00)"
$ws.Range("F4").Value = "Автобус 3 (08:00 - 16:00)"
$ws.Range("G4").Value = "Автобус 4 (10:00 - 18:00)"

$ws.Range("C5").Value = "Автобус 8 (09:00 - 17:00)"
$ws.Range("D5").Value = "Автобус 8 (06:00 - 14:00)"
$ws.Range("E5").Value = "Автобус 1 (09:00 - 17:00)"
$ws.Range("F5").Value = "Автобус 7 (09:00 - 17:00)"
$ws.Range("G5").Value = "Автобус 5 (06:00 - 14:00)"

$ws.Range("C6").Value = "Автобус 5 (10:00 - 18:00)"
$ws.Range("D6").Value = "Автобус 7 (07:00 - 15:00)"
$ws.Range("E6").Value = "Автобус 4 (08:00 - 16:00)"
$ws.Range("F6").Value = "Автобус 1 (10:00 - 18:00)"
$ws.Range("G6").Value = "Автобус 2 (08:00 - 16:00)"

$ws.Range("C7").Value = "Автобус 2 (08:00 - 16:00)"
$ws.Range("D7").Value = "Автобус 4 (10:00 - 18:00)"
$ws.Range("E7").Value = "Автобус 4 (06:00 - 14:00)"
$ws.Range("F7").Value = "Автобус 8 (09:00 - 17:00)"
$ws.Range("G7").Value = "Автобус 8 (08:00 - 16:00)"

$ws.Range("C8").Value = "Автобус 5 (10:00 - 22:00)"
$ws.Range("F8").Value = "Автобус 2 (09:00 - 21:00)"
$ws.Range("I8").Value = "Автобус 1 (16:00 - 04:00)"

$ws.Range("C9").Value = "Автобус 8 (09:00 - 21:00)"
$ws.Range("I9").Value = "Автобус 1 (05:00 - 17:00)"

$ws.Range("D10").Value = "Автобус 5 (06:00 - 18:00)"
$ws.Range("G10").Value = "Автобус 7 (04:00 - 16:00)"

$ws.Range("D11").Value = "Автобус 6 (10:00 - 22:00)"
$ws.Range("G11").Value = "Автобус 2 (07:00 - 19:00)"

$ws.Range("E12").Value = "Автобус 8 (16:00 - 04:00)"
$ws.Range("H12").Value = "Автобус 7 (15:00 - 03:00)"

$ws.Range("E13").Value = "Автобус 6 (01:00 - 13:00)"
$ws.Range("H13").Value = "Автобус 4 (05:00 - 17:00)"
